# Updates the "Price" (column D) and "Volume(1h)" (column E) figures on the
# cryptos sheet to a newer snapshot. Values are written with NumberFormat
# forced to Text ("@") so figures like "66.531.20" / "0.999" / "9.20" are
# preserved as literal text (matching the original cell contents, which are
# themselves plain text) instead of being auto-converted/rounded as numbers.
# The style is then reset back to "Normal" so we don't leave a stray
# text-number-format style applied to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '66.531.20'
$cell.Style = 'Normal'

$cell = $ws.Range('E2')
$cell.NumberFormat = '@'
$cell.Value = '  +1.04%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '3.302.06'
$cell.Style = 'Normal'

$cell = $ws.Range('E3')
$cell.NumberFormat = '@'
$cell.Value = '  +0.25%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E4')
$cell.NumberFormat = '@'
$cell.Value = '  -0.05%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '588.57'
$cell.Style = 'Normal'

$cell = $ws.Range('E5')
$cell.NumberFormat = '@'
$cell.Value = '  +2.65%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '180.88'
$cell.Style = 'Normal'

$cell = $ws.Range('E6')
$cell.NumberFormat = '@'
$cell.Value = '  +1.33%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.639'
$cell.Style = 'Normal'

$cell = $ws.Range('E7')
$cell.NumberFormat = '@'
$cell.Value = '  +1.57%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E8')
$cell.NumberFormat = '@'
$cell.Value = '  -0.05%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '3.295.79'
$cell.Style = 'Normal'

$cell = $ws.Range('E9')
$cell.NumberFormat = '@'
$cell.Value = '  +0.13%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E10')
$cell.NumberFormat = '@'
$cell.Value = '  +0.18%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E11')
$cell.NumberFormat = '@'
$cell.Value = '  +2.86%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '0.403'
$cell.Style = 'Normal'

$cell = $ws.Range('E12')
$cell.NumberFormat = '@'
$cell.Value = '  +0.25%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '3.876.92'
$cell.Style = 'Normal'

$cell = $ws.Range('E13')
$cell.NumberFormat = '@'
$cell.Value = '  +0.21%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E14')
$cell.NumberFormat = '@'
$cell.Value = '  -2.20%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '66.522.00'
$cell.Style = 'Normal'

$cell = $ws.Range('E15')
$cell.NumberFormat = '@'
$cell.Value = '  +0.75%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '26.67'
$cell.Style = 'Normal'

$cell = $ws.Range('E16')
$cell.NumberFormat = '@'
$cell.Value = '  +0.34%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E17')
$cell.NumberFormat = '@'
$cell.Value = '  -0.13%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '3.275.50'
$cell.Style = 'Normal'

$cell = $ws.Range('E18')
$cell.NumberFormat = '@'
$cell.Value = '  -1.90%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '428.69'
$cell.Style = 'Normal'

$cell = $ws.Range('E19')
$cell.NumberFormat = '@'
$cell.Value = '  -1.32%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '5.47'
$cell.Style = 'Normal'

$cell = $ws.Range('E20')
$cell.NumberFormat = '@'
$cell.Value = '  -2.29%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '13.04'
$cell.Style = 'Normal'

$cell = $ws.Range('E21')
$cell.NumberFormat = '@'
$cell.Value = '  -2.11%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '7.29'
$cell.Style = 'Normal'

$cell = $ws.Range('E22')
$cell.NumberFormat = '@'
$cell.Value = '  -1.88%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E23')
$cell.NumberFormat = '@'
$cell.Value = '  +0.17%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '71.51'
$cell.Style = 'Normal'

$cell = $ws.Range('E24')
$cell.NumberFormat = '@'
$cell.Value = '  -1.32%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '5.74'
$cell.Style = 'Normal'

$cell = $ws.Range('E25')
$cell.NumberFormat = '@'
$cell.Value = '  +0.95%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '3.448.85'
$cell.Style = 'Normal'

$cell = $ws.Range('E26')
$cell.NumberFormat = '@'
$cell.Value = '  +0.01%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E27')
$cell.NumberFormat = '@'
$cell.Value = '  +0.03%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E28')
$cell.NumberFormat = '@'
$cell.Value = '  +5.35%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '0.0000115'
$cell.Style = 'Normal'

$cell = $ws.Range('E29')
$cell.NumberFormat = '@'
$cell.Value = '  +0.44%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '9.20'
$cell.Style = 'Normal'

$cell = $ws.Range('E30')
$cell.NumberFormat = '@'
$cell.Value = '  +2.53%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'

$cell = $ws.Range('E31')
$cell.NumberFormat = '@'
$cell.Value = '  -0.34%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '1.92'
$cell.Style = 'Normal'

$cell = $ws.Range('E32')
$cell.NumberFormat = '@'
$cell.Value = '  -1.27%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '22.34'
$cell.Style = 'Normal'

$cell = $ws.Range('E33')
$cell.NumberFormat = '@'
$cell.Value = '  -0.45%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E34')
$cell.NumberFormat = '@'
$cell.Value = '  +0.08%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '5.16'
$cell.Style = 'Normal'

$cell = $ws.Range('E35')
$cell.NumberFormat = '@'
$cell.Value = '  +0.14%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '6.56'
$cell.Style = 'Normal'

$cell = $ws.Range('E36')
$cell.NumberFormat = '@'
$cell.Value = '  -1.23%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E37')
$cell.NumberFormat = '@'
$cell.Value = '  -1.02%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '158.81'
$cell.Style = 'Normal'

$cell = $ws.Range('E38')
$cell.NumberFormat = '@'
$cell.Value = '  +1.08%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E39')
$cell.NumberFormat = '@'
$cell.Value = '  -1.74%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '2.862.15'
$cell.Style = 'Normal'

$cell = $ws.Range('E40')
$cell.NumberFormat = '@'
$cell.Value = '  +2.97%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E41')
$cell.NumberFormat = '@'
$cell.Value = '  +0.09%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '26.31'
$cell.Style = 'Normal'

$cell = $ws.Range('E42')
$cell.NumberFormat = '@'
$cell.Value = '  -2.93%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '4.32'
$cell.Style = 'Normal'

$cell = $ws.Range('E43')
$cell.NumberFormat = '@'
$cell.Value = '  -0.71%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E44')
$cell.NumberFormat = '@'
$cell.Value = '  -4.25%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E45')
$cell.NumberFormat = '@'
$cell.Value = '  -1.70%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '0.0656'
$cell.Style = 'Normal'

$cell = $ws.Range('E46')
$cell.NumberFormat = '@'
$cell.Value = '  -0.47%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '5.93'
$cell.Style = 'Normal'

$cell = $ws.Range('E47')
$cell.NumberFormat = '@'
$cell.Value = '  -3.05%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '2.31'
$cell.Style = 'Normal'

$cell = $ws.Range('E48')
$cell.NumberFormat = '@'
$cell.Value = '  +0.46%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '313.18'
$cell.Style = 'Normal'

$cell = $ws.Range('E49')
$cell.NumberFormat = '@'
$cell.Value = '  -3.07%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '22.85'
$cell.Style = 'Normal'

$cell = $ws.Range('E50')
$cell.NumberFormat = '@'
$cell.Value = '  -3.15%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E51')
$cell.NumberFormat = '@'
$cell.Value = '  -0.11%  '
$cell.Style = 'Normal'

